# updated to check all misc pages for IE11 message, even if one is missing
#
# The "ComparePages" sheet listed two Lexus series that are no longer sold
# (GS and GS F) among the car-series compare URLs. Remove those two rows so
# the sheet now only lists the currently reviewed series, then leave the
# sheet active/selected the way the author left Excel (on ComparePages,
# with the row that used to hold the LC-hybrid/"GSF" boundary selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ComparePages")

# Row 4 holds ".../compare/series/GS" and row 14 holds
# ".../compare/series/GSF" (discontinued models) - delete both. Delete the
# later row first so the row 4 index used for the first delete is still
# valid.
$ws.Rows("14").Delete()
$ws.Rows("4").Delete()

# Leave the workbook with ComparePages as the active/selected sheet, with
# the former "GSF" row position selected.
$ws.Activate()
$ws.Range("A13:XFD13").Select()
